$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts Volume/P_E/EPS(TTM) right by one
# and leaves room for "Market Cap")
$ws.Range("D1:D2").EntireColumn.Insert()

# New column D: Market Cap
$ws.Range("D1").Value = "Market Cap"
$ws.Range("D2").Value = '''$302.66B'

# New column H: Div. Yield (after the insert, old F/G live at G/H already;
# H is brand new)
$ws.Range("H1").Value = "Div. Yield"
$ws.Range("H2").Value = '''2.88%'

# Updated price / change figures for the existing row
$ws.Range("B2").Value = '''$70.37'
$ws.Range("C2").Value = "-0.37(0.52%) 1D"

# Updated P_E (shifted from E to F by the column insert)
$ws.Range("F2").Value = '''28.5'

# Match the bold/centered/bordered header style used by the other header cells
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
